$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = -20.926
$ws.Range("A12").Value = -21.694
$ws.Range("E13").Value = 12.817
$ws.Range("A18").Value = -21.694
